$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Two new inserts that replaced the old "insert 2 pokemon from Generation 3" entry
$ws.Range("B28").Value = "Insert Treecko"
$ws.Range("B29").Value = "Insert Mudkip"
$ws.Range("A28").Value = $true
$ws.Range("A29").Value = $true

# Two new deletes that replaced the old "(possibly) delete the 2 new insertions from Gen 3" entry
$ws.Range("B25").Value = "Delete Treecko"
$ws.Range("B26").Value = "Delete Mudkip"
$ws.Range("A25").Value = $true
$ws.Range("A26").Value = $true
